$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, and whether the text
# looks like a plain number (needs a quote-prefix so Excel keeps it
# as text instead of silently converting it to a numeric value).
$updates = @(
    @{ Cell = 'D2'; Value = '29.311.70'; Numeric = $false },
    @{ Cell = 'E2'; Value = '  +0.00%  '; Numeric = $false },
    @{ Cell = 'D3'; Value = '1.875.83'; Numeric = $false },
    @{ Cell = 'E3'; Value = '  +0.13%  '; Numeric = $false },
    @{ Cell = 'E4'; Value = '  +0.02%  '; Numeric = $false },
    @{ Cell = 'D5'; Value = '0.7125'; Numeric = $true },
    @{ Cell = 'E5'; Value = '  +0.13%  '; Numeric = $false },
    @{ Cell = 'D6'; Value = '242.46'; Numeric = $true },
    @{ Cell = 'E6'; Value = '  +0.33%  '; Numeric = $false },
    @{ Cell = 'E7'; Value = '  +0.04%  '; Numeric = $false },
    @{ Cell = 'D8'; Value = '0.08023'; Numeric = $true },
    @{ Cell = 'E8'; Value = '  +3.22%  '; Numeric = $false },
    @{ Cell = 'D9'; Value = '0.3153'; Numeric = $true },
    @{ Cell = 'E9'; Value = '  +1.43%  '; Numeric = $false },
    @{ Cell = 'D10'; Value = '24.99'; Numeric = $true },
    @{ Cell = 'E10'; Value = '  -0.57%  '; Numeric = $false },
    @{ Cell = 'D11'; Value = '0.08228'; Numeric = $true },
    @{ Cell = 'E11'; Value = '  -2.00%  '; Numeric = $false },
    @{ Cell = 'D12'; Value = '1.880.91'; Numeric = $false },
    @{ Cell = 'E12'; Value = '  +0.62%  '; Numeric = $false },
    @{ Cell = 'B13'; Value = 'Litecoin'; Numeric = $false },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; Numeric = $false },
    @{ Cell = 'D13'; Value = '94.84'; Numeric = $true },
    @{ Cell = 'E13'; Value = '  +4.04%  '; Numeric = $false },
    @{ Cell = 'B14'; Value = 'Polkadot'; Numeric = $false },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; Numeric = $false },
    @{ Cell = 'D14'; Value = '5.246'; Numeric = $true },
    @{ Cell = 'E14'; Value = '  +0.12%  '; Numeric = $false },
    @{ Cell = 'D15'; Value = '0.7121'; Numeric = $true },
    @{ Cell = 'E15'; Value = '  +0.11%  '; Numeric = $false },
    @{ Cell = 'D16'; Value = '6.412'; Numeric = $true },
    @{ Cell = 'D17'; Value = '0.000008523'; Numeric = $true },
    @{ Cell = 'E17'; Value = '  +3.96%  '; Numeric = $false },
    @{ Cell = 'D18'; Value = '29.327.63'; Numeric = $false },
    @{ Cell = 'E18'; Value = '  +0.01%  '; Numeric = $false },
    @{ Cell = 'D19'; Value = '243.71'; Numeric = $true },
    @{ Cell = 'E19'; Value = '  +1.50%  '; Numeric = $false },
    @{ Cell = 'D20'; Value = '2.142.12'; Numeric = $false },
    @{ Cell = 'E20'; Value = '  +1.01%  '; Numeric = $false },
    @{ Cell = 'D21'; Value = '13.23'; Numeric = $true },
    @{ Cell = 'E21'; Value = '  +0.25%  '; Numeric = $false },
    @{ Cell = 'E22'; Value = '  +0.02%  '; Numeric = $false },
    @{ Cell = 'D23'; Value = '7.774'; Numeric = $true },
    @{ Cell = 'E23'; Value = '  +0.20%  '; Numeric = $false },
    @{ Cell = 'E24'; Value = '  +0.00%  '; Numeric = $false },
    @{ Cell = 'E25'; Value = '  -1.57%  '; Numeric = $false },
    @{ Cell = 'B26'; Value = 'Monero'; Numeric = $false },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; Numeric = $false },
    @{ Cell = 'D26'; Value = '162.47'; Numeric = $true },
    @{ Cell = 'E26'; Value = '  -0.36%  '; Numeric = $false },
    @{ Cell = 'B27'; Value = 'Cosmos'; Numeric = $false },
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; Numeric = $false },
    @{ Cell = 'D27'; Value = '9.037'; Numeric = $true },
    @{ Cell = 'E27'; Value = '  +0.15%  '; Numeric = $false },
    @{ Cell = 'D28'; Value = '18.50'; Numeric = $true },
    @{ Cell = 'E28'; Value = '  +0.02%  '; Numeric = $false },
    @{ Cell = 'D29'; Value = '1.502'; Numeric = $true },
    @{ Cell = 'E29'; Value = '  -0.46%  '; Numeric = $false },
    @{ Cell = 'D30'; Value = '4.411'; Numeric = $true },
    @{ Cell = 'E30'; Value = '  +0.14%  '; Numeric = $false },
    @{ Cell = 'D31'; Value = '4.308'; Numeric = $true },
    @{ Cell = 'E31'; Value = '  -0.28%  '; Numeric = $false },
    @{ Cell = 'D32'; Value = '0.05368'; Numeric = $true },
    @{ Cell = 'E32'; Value = '  +1.39%  '; Numeric = $false },
    @{ Cell = 'D33'; Value = '1.171'; Numeric = $true },
    @{ Cell = 'E33'; Value = '  -9.13%  '; Numeric = $false },
    @{ Cell = 'D34'; Value = '1.942'; Numeric = $true },
    @{ Cell = 'E34'; Value = '  +0.14%  '; Numeric = $false },
    @{ Cell = 'D35'; Value = '0.7631'; Numeric = $true },
    @{ Cell = 'E35'; Value = '  +2.57%  '; Numeric = $false },
    @{ Cell = 'D36'; Value = '1.178'; Numeric = $true },
    @{ Cell = 'E36'; Value = '  +0.06%  '; Numeric = $false },
    @{ Cell = 'D37'; Value = '2.689'; Numeric = $true },
    @{ Cell = 'E37'; Value = '  -0.54%  '; Numeric = $false },
    @{ Cell = 'D38'; Value = '0.01876'; Numeric = $true },
    @{ Cell = 'E38'; Value = '  -0.02%  '; Numeric = $false },
    @{ Cell = 'D39'; Value = '1.265.04'; Numeric = $false },
    @{ Cell = 'E39'; Value = '  +3.07%  '; Numeric = $false },
    @{ Cell = 'D40'; Value = '2.754'; Numeric = $true },
    @{ Cell = 'E40'; Value = '  +0.93%  '; Numeric = $false },
    @{ Cell = 'D41'; Value = '6.475'; Numeric = $true },
    @{ Cell = 'E41'; Value = '  -0.93%  '; Numeric = $false },
    @{ Cell = 'D42'; Value = '0.9139'; Numeric = $true },
    @{ Cell = 'E42'; Value = '  +3.28%  '; Numeric = $false },
    @{ Cell = 'D43'; Value = '112.66'; Numeric = $true },
    @{ Cell = 'E43'; Value = '  +3.21%  '; Numeric = $false },
    @{ Cell = 'D44'; Value = '74.06'; Numeric = $true },
    @{ Cell = 'E44'; Value = '  +2.23%  '; Numeric = $false },
    @{ Cell = 'D45'; Value = '0.00000000134'; Numeric = $true },
    @{ Cell = 'E45'; Value = '  +9.49%  '; Numeric = $false },
    @{ Cell = 'E46'; Value = '  +0.03%  '; Numeric = $false },
    @{ Cell = 'D47'; Value = '2.038.62'; Numeric = $false },
    @{ Cell = 'E47'; Value = '  +1.06%  '; Numeric = $false },
    @{ Cell = 'D48'; Value = '0.5229'; Numeric = $true },
    @{ Cell = 'E48'; Value = '  +0.69%  '; Numeric = $false },
    @{ Cell = 'D49'; Value = '1.799'; Numeric = $true },
    @{ Cell = 'E49'; Value = '  +0.06%  '; Numeric = $false },
    @{ Cell = 'D50'; Value = '9.475'; Numeric = $true },
    @{ Cell = 'E50'; Value = '  +1.04%  '; Numeric = $false },
    @{ Cell = 'D51'; Value = '0.4353'; Numeric = $true },
    @{ Cell = 'E51'; Value = '  +1.02%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # Leading apostrophe forces text storage; reset the style
        # afterwards so no stray number-format/quote-prefix style
        # is left behind on the cell.
        $rng.Value = "'" + $u.Value
        $rng.Style = 'Normal'
    } else {
        $rng.Value = $u.Value
    }
}
